# Added some poisson loo numbers
# Fill in the LOO-CV results for the two newly-run models:
#   row 33 -> 09p_model_fit_constant_zeta_full_x
#   row 35 -> 09r_model_fit_no_zeta_full_x

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 33: 09p_model_fit_constant_zeta_full_x ------------------------
# "Good k values" is low here (0.67), so Excel highlights it with a
# solid dark-red fill (matching the existing warning convention in the sheet).
$ws.Range("D33").Value = 0.67
$ws.Range("D33").NumberFormat = "0.00%"
$ws.Range("D33").Interior.Color = 192

$ws.Range("E33").Value = -354064
$ws.Range("E33").NumberFormat = "#,##0"

$ws.Range("F33").Value = 48022
$ws.Range("F33").NumberFormat = "#,##0"

$ws.Range("G33").Value = 708129
$ws.Range("G33").NumberFormat = "#,##0"

# ---- Row 35: 09r_model_fit_no_zeta_full_x -------------------------------
$ws.Range("D35").Value = 0.967
$ws.Range("D35").NumberFormat = "0.00%"

$ws.Range("E35").Value = -355560
$ws.Range("E35").NumberFormat = "#,##0"

$ws.Range("F35").Value = 14532
$ws.Range("F35").NumberFormat = "#,##0"

$ws.Range("G35").Value = 711120
$ws.Range("G35").NumberFormat = "#,##0"

# H35 already carries the "#,##0.00"-with-bottom-border style used throughout
# this block; it was previously blank and now records the susc_scaling value.
$ws.Range("H35").Value = 0.95

# Move the selection/scroll position to where the new numbers were entered.
$ws.Range("D53").Select() | Out-Null
